$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 116
$ws.Range("I2").Value = 367
$ws.Range("J2").Value = 1386
$ws.Range("K2").Value = 7
$ws.Range("L2").Value = 357
$ws.Range("M2").Value = 21
$ws.Range("N2").Value = 234
$ws.Range("P2").Value = 4
$ws.Range("R2").Value = 21
$ws.Range("S2").Value = 153
$ws.Range("T2").Value = 232
$ws.Range("V2").Value = 2101
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 2105
$ws.Range("Y2").Value = 5
$ws.Range("Z2").Value = 27
$ws.Range("AA2").Value = 16
